# Auto-generated edit script: updates cached numeric values in several
# worksheets per the scheduled-runner diff (market-price / profit refresh).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 199
$ws.Range("I4").Value = 199
$ws.Range("K4").Value = 199
$ws.Range("M4").Value = -85
$ws.Range("H9").Value = 30.333334
$ws.Range("I9").Value = 30.5
$ws.Range("K9").Value = 30.5
$ws.Range("M9").Value = 138.5
$ws.Range("H33").Value = 765.1539
$ws.Range("I33").Value = 816.4167
$ws.Range("K33").Value = 816.4167
$ws.Range("M33").Value = -587.4167
$ws.Range("H40").Value = 2267
$ws.Range("I40").Value = 1950.5
$ws.Range("K40").Value = 1950.5
$ws.Range("M40").Value = -1775.5
$ws.Range("H98").Value = 6521.75
$ws.Range("I98").Value = 1540.1538
$ws.Range("K98").Value = 1540.1538
$ws.Range("M98").Value = -42.15380000000005
$ws.Range("H122").Value = 6521.75
$ws.Range("I122").Value = 1540.1538
$ws.Range("K122").Value = 4620.4614
$ws.Range("M122").Value = -2170.4614
$ws.Range("H137").Value = 7165.0415
$ws.Range("I137").Value = 6791.9414
$ws.Range("K137").Value = 20375.8242
$ws.Range("M137").Value = -17825.8242
$ws.Range("H138").Value = 7000
$ws.Range("I138").Value = 3000
$ws.Range("K138").Value = 9000
$ws.Range("M138").Value = -3860

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 6000
$ws.Range("J23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("N23").Value = -6518
$ws.Range("H61").Value = 8599.799999999999
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20424
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H132").Value = 8652.200000000001
$ws.Range("I132").Value = 1087
$ws.Range("K132").Value = 3261
$ws.Range("M132").Value = -731
$ws.Range("H136").Value = 8599.799999999999
$ws.Range("J136").Value = 20000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -65100
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200
$ws.Range("H139").Value = 74205
$ws.Range("J139").Value = 74205
$ws.Range("L139").Value = 74205
$ws.Range("N139").Value = -84485

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1550
$ws.Range("I16").Value = 1550
$ws.Range("K16").Value = 1550
$ws.Range("M16").Value = -1263
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 7286.5713
$ws.Range("I31").Value = 5681.1
$ws.Range("J31").Value = 11300.25
$ws.Range("K31").Value = 5681.1
$ws.Range("L31").Value = 11300.25
$ws.Range("M31").Value = -5386.1
$ws.Range("N31").Value = -11890.25
$ws.Range("H34").Value = 7286.5713
$ws.Range("I34").Value = 5681.1
$ws.Range("J34").Value = 11300.25
$ws.Range("K34").Value = 5681.1
$ws.Range("L34").Value = 11300.25
$ws.Range("M34").Value = -5479.1
$ws.Range("N34").Value = -11704.25
$ws.Range("H51").Value = 90
$ws.Range("I51").Value = 90
$ws.Range("K51").Value = 90
$ws.Range("M51").Value = 646
$ws.Range("H59").Value = 38248.75
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -47290
$ws.Range("H61").Value = 90
$ws.Range("I61").Value = 90
$ws.Range("K61").Value = 90
$ws.Range("M61").Value = 258
$ws.Range("H113").Value = 1550
$ws.Range("I113").Value = 1550
$ws.Range("K113").Value = 1550
$ws.Range("M113").Value = 620
$ws.Range("H134").Value = 7034.7
$ws.Range("I134").Value = 4543.375
$ws.Range("K134").Value = 13630.125
$ws.Range("M134").Value = -11095.125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 198.33333
$ws.Range("J23").Value = 99.5
$ws.Range("L23").Value = 298.5
$ws.Range("N23").Value = -768.5
$ws.Range("H24").Value = 350
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = -370
$ws.Range("N24").Value = -1960
$ws.Range("H34").Value = 848.75
$ws.Range("I34").Value = 131.66667
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 395.00001
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -311.00001
$ws.Range("N34").Value = -9168
$ws.Range("H38").Value = 138.71428
$ws.Range("J38").Value = 160
$ws.Range("L38").Value = 480
$ws.Range("N38").Value = -1174
$ws.Range("H50").Value = 238.33333
$ws.Range("I50").Value = 143.125
$ws.Range("J50").Value = 1000
$ws.Range("K50").Value = 429.375
$ws.Range("L50").Value = 3000
$ws.Range("M50").Value = 51.625
$ws.Range("N50").Value = -3962
$ws.Range("H53").Value = 238.33333
$ws.Range("I53").Value = 143.125
$ws.Range("J53").Value = 1000
$ws.Range("K53").Value = 429.375
$ws.Range("L53").Value = 3000
$ws.Range("M53").Value = 51.625
$ws.Range("N53").Value = -3962
$ws.Range("H68").Value = 940.6667
$ws.Range("I68").Value = 940.6667
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2822.0001
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2011.0001
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 940.6667
$ws.Range("I71").Value = 940.6667
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8466.0003
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4410.0003
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 1144.5
$ws.Range("I80").Value = 1144.5
$ws.Range("K80").Value = 3433.5
$ws.Range("M80").Value = -2497.5
$ws.Range("H83").Value = 1144.5
$ws.Range("I83").Value = 1144.5
$ws.Range("K83").Value = 10300.5
$ws.Range("M83").Value = -5620.5
$ws.Range("H121").Value = 1012.75
$ws.Range("I121").Value = 245
$ws.Range("J121").Value = 1268.6666
$ws.Range("K121").Value = 735
$ws.Range("L121").Value = 3805.9998
$ws.Range("M121").Value = 575
$ws.Range("N121").Value = -6425.9998
$ws.Range("H132").Value = 1617.5
$ws.Range("I132").Value = 1485
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 13365
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -10835
$ws.Range("N132").Value = -20810

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3320.8333
$ws.Range("I102").Value = 2983.3333
$ws.Range("K102").Value = 2983.3333
$ws.Range("M102").Value = -1361.3333

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H46").Value = 6552.9414
$ws.Range("I46").Value = 7450
$ws.Range("K46").Value = 7450
$ws.Range("M46").Value = -7262
$ws.Range("H132").Value = 8096.533
$ws.Range("I132").Value = 6745.6665
$ws.Range("K132").Value = 20236.9995
$ws.Range("M132").Value = -17706.9995
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H137").Value = 30000
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -24900
$ws.Range("N137").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 15000400
$ws.Range("I19").Value = 800
$ws.Range("K19").Value = 800
$ws.Range("M19").Value = -626
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H81").Value = 1067.8
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1067.8
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1770.6154
$ws.Range("I122").Value = 1328.6
$ws.Range("J122").Value = 2046.875
$ws.Range("K122").Value = 3985.8
$ws.Range("L122").Value = 6140.625
$ws.Range("M122").Value = -1535.8
$ws.Range("N122").Value = -11040.625
$ws.Range("H132").Value = 11372.5
$ws.Range("I132").Value = 4745
$ws.Range("K132").Value = 14235
$ws.Range("M132").Value = -11705
